# Northern Ireland changed its wording again -- refresh the UK daily
# case/death figures on Sheet1 and highlight the revised NI column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily snapshot values for row 2 (DateVal, TotalUKCases, NewUKCases,
# TotalUKDeaths, EnglandCases, ScotlandCases, WalesCases, NICases).
$ws.Range("A2").Value = 43909
$ws.Range("B2").Value = 3269
$ws.Range("C2").Value = 643
$ws.Range("D2").Value = 144
$ws.Range("E2").Value = 2756
$ws.Range("F2").Value = 266
$ws.Range("G2").Value = 170
$ws.Range("H2").Value = 77

# NICases (D2) picks up a highlight fill now that NI revised its reporting
# again, same light background used for the header row.
$ws.Range("D2").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1

# Move the active selection.
$ws.Range("B5").Select()
